# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell that still
# showed the old "Ready for handoff" status is refreshed to "In
# Translation", and the Status / per-language columns that held that text
# are re-sized to fit the (now shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells -------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsOverview.Columns("E:F").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (col C) ----------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsZhCn.Columns("C:C").ColumnWidth = 12.5

# --- de-de sheet: Status column (col C) ----------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
